$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original header texts before overwriting anything.
$hA = $ws.Range("A1").Text
$hB = $ws.Range("B1").Text
$hC = $ws.Range("C1").Text
$hD = $ws.Range("D1").Text
$hE = $ws.Range("E1").Text
$hF = $ws.Range("F1").Text

# Shift headers one column to the right; new first column is the monthly
# income header that used to live in G1.
$ws.Range("G1").Value = $hF
$ws.Range("F1").Value = $hE
$ws.Range("E1").Value = $hD
$ws.Range("D1").Value = $hC
$ws.Range("C1").Value = $hB
$ws.Range("B1").Value = $hA
$ws.Range("A1").Value = "Rendimento Mensal (R$)"

# New row 2 dataset
$ws.Range("A2").Value = 10000
$ws.Range("B2").Value = 120000
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 7000
$ws.Range("G2").Value = 20196

# Remove old extra data rows 3-8
$ws.Range("A3:G8").EntireRow.Delete()
